# Insert a new row for "Nason Creek Lower 11" above the current row 61
# ("Nason Creek Lower 12"), shifting all subsequent rows down by one.
# This corresponds to integrating new Okanogan EDT results into the
# Steelhead Habitat Quality (Restoration) worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 61; existing rows 61-102 shift to 62-103.
$ws.Rows("61:61").Insert()

# Populate the newly inserted row with the new reach's data.
$ws.Range("A61").Value = "Nason Creek Lower 11"
$ws.Range("B61").Value = "Wenatchee"
$ws.Range("C61").Value = "Lower Nason Creek"
$ws.Range("D61").Value = "yes"
$ws.Range("E61").Value = "yes"
$ws.Range("F61").Value = "yes"
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 3
$ws.Range("I61").Value = 3
$ws.Range("J61").Value = 5
$ws.Range("K61").Value = 3
$ws.Range("L61").Value = 3
$ws.Range("M61").Value = 3
$ws.Range("N61").Value = 1
$ws.Range("O61").Value = 1
$ws.Range("P61").Value = 1
$ws.Range("Q61").Value = 3
$ws.Range("R61").Value = 2
$ws.Range("S61").Value = 1
$ws.Range("T61").Value = 22
$ws.Range("U61").Value = 0.4888888888888889
$ws.Range("V61").Value = 5
$ws.Range("W61").Value = 1
$ws.Range("X61").Value = "Off-Channel-Side-Channels,PoolQuantity&Quality,Temperature-Rearing"
$ws.Range("Y61").Value = "Stability,Cover-Wood,Flow-SummerBaseFlow,Off-Channel-Floodplain,Riparian"
$ws.Range("Z61").Value = "Stability,Cover-Wood,Flow-SummerBaseFlow,Off-Channel-Floodplain,Off-Channel-Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing"
